$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1140.4524
$ws.Range("J17").Value = 1148.7711
$ws.Range("L17").Value = 3446.3133
$ws.Range("N17").Value = -3782.3133
$ws.Range("H28").Value = 775.5455
$ws.Range("I28").Value = 797.4211
$ws.Range("J28").Value = 745.8571
$ws.Range("K28").Value = 797.4211
$ws.Range("L28").Value = 745.8571
$ws.Range("M28").Value = -312.4211
$ws.Range("N28").Value = -1715.8571
$ws.Range("H70").Value = 3072.7273
$ws.Range("I70").Value = 1600
$ws.Range("J70").Value = 3625
$ws.Range("K70").Value = 4800
$ws.Range("L70").Value = 10875
$ws.Range("M70").Value = -4530
$ws.Range("N70").Value = -11415
$ws.Range("H73").Value = 3072.7273
$ws.Range("I73").Value = 1600
$ws.Range("J73").Value = 3625
$ws.Range("K73").Value = 4800
$ws.Range("L73").Value = 10875
$ws.Range("M73").Value = -3864
$ws.Range("N73").Value = -12747
$ws.Range("H101").Value = 2278.1
$ws.Range("I101").Value = 245
$ws.Range("J101").Value = 3633.5
$ws.Range("K101").Value = 735
$ws.Range("L101").Value = 10900.5
$ws.Range("M101").Value = 887
$ws.Range("N101").Value = -14144.5
$ws.Range("H115").Value = 1077.9166
$ws.Range("I115").Value = 637.2222
$ws.Range("J115").Value = 2400
$ws.Range("K115").Value = 1911.6666
$ws.Range("L115").Value = 7200
$ws.Range("M115").Value = -344.6666
$ws.Range("N115").Value = -10334
$ws.Range("H116").Value = 1400
$ws.Range("I116").Value = 1700
$ws.Range("J116").Value = 1100
$ws.Range("K116").Value = 1700
$ws.Range("L116").Value = 1100
$ws.Range("M116").Value = 1742
$ws.Range("N116").Value = -7984
$ws.Range("H127").Value = 859
$ws.Range("I127").Value = 418.0909
$ws.Range("J127").Value = 1162.125
$ws.Range("K127").Value = 1254.2727
$ws.Range("L127").Value = 3486.375
$ws.Range("M127").Value = 3705.7273
$ws.Range("N127").Value = -13406.375
$ws.Range("H129").Value = 1000
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 1000
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 3000
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -13000
$ws.Range("H137").Value = 22728416
$ws.Range("I137").Value = 1123.2941
$ws.Range("J137").Value = 100001210
$ws.Range("K137").Value = 3369.8823
$ws.Range("L137").Value = 300003630
$ws.Range("M137").Value = -819.8823000000002
$ws.Range("N137").Value = -300008730
$ws.Range("H138").Value = 1401.091
$ws.Range("I138").Value = 1243.2463
$ws.Range("J138").Value = 2762.5
$ws.Range("K138").Value = 3729.7389
$ws.Range("L138").Value = 8287.5
$ws.Range("M138").Value = 1410.2611
$ws.Range("N138").Value = -18567.5

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8637.242
$ws.Range("I32").Value = 4921.316
$ws.Range("J32").Value = 27464.6
$ws.Range("K32").Value = 4921.316
$ws.Range("L32").Value = 27464.6
$ws.Range("M32").Value = -4634.316
$ws.Range("N32").Value = -28038.6
$ws.Range("H68").Value = 59952.6
$ws.Range("J68").Value = 59952.6
$ws.Range("L68").Value = 59952.6
$ws.Range("N68").Value = -61574.6
$ws.Range("H71").Value = 59952.6
$ws.Range("J71").Value = 59952.6
$ws.Range("L71").Value = 179857.8
$ws.Range("N71").Value = -187969.8
$ws.Range("H74").Value = 5145.7646
$ws.Range("I74").Value = 955.9643
$ws.Range("J74").Value = 24698.166
$ws.Range("K74").Value = 955.9643
$ws.Range("L74").Value = 24698.166
$ws.Range("M74").Value = -81.96429999999998
$ws.Range("N74").Value = -26446.166
$ws.Range("H77").Value = 5145.7646
$ws.Range("I77").Value = 955.9643
$ws.Range("J77").Value = 24698.166
$ws.Range("K77").Value = 4779.8215
$ws.Range("L77").Value = 123490.83
$ws.Range("M77").Value = -411.8215
$ws.Range("N77").Value = -132226.83
$ws.Range("H132").Value = 20695.076
$ws.Range("I132").Value = 1387.8529
$ws.Range("J132").Value = 57164.277
$ws.Range("K132").Value = 4163.5587
$ws.Range("L132").Value = 171492.831
$ws.Range("M132").Value = -1633.5587
$ws.Range("N132").Value = -176552.831

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 20330.254
$ws.Range("I134").Value = 22497.303
$ws.Range("J134").Value = 1188
$ws.Range("K134").Value = 67491.909
$ws.Range("L134").Value = 3564
$ws.Range("M134").Value = -64956.909
$ws.Range("N134").Value = -8634

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4467546
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("H34").Value = 4467546
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H132").Value = 1056.8414
$ws.Range("I132").Value = 803.8095
$ws.Range("K132").Value = 2411.4285
$ws.Range("M132").Value = 118.5715

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 898.67346
$ws.Range("I5").Value = 473.08334
$ws.Range("J5").Value = 1036.7028
$ws.Range("K5").Value = 1419.25002
$ws.Range("L5").Value = 3110.1084
$ws.Range("M5").Value = -1307.25002
$ws.Range("N5").Value = -3334.1084
$ws.Range("H122").Value = 11906101
$ws.Range("I122").Value = 15152126
$ws.Range("J122").Value = 4008.3333
$ws.Range("K122").Value = 136369134
$ws.Range("L122").Value = 36074.9997
$ws.Range("M122").Value = -136366684
$ws.Range("N122").Value = -40974.9997
$ws.Range("H131").Value = 3433.8538
$ws.Range("I131").Value = 4452.6665
$ws.Range("J131").Value = 2846.077
$ws.Range("K131").Value = 13357.9995
$ws.Range("L131").Value = 8538.231
$ws.Range("M131").Value = -8317.999500000002
$ws.Range("N131").Value = -18618.231
$ws.Range("H135").Value = 898.67346
$ws.Range("I135").Value = 473.08334
$ws.Range("J135").Value = 1036.7028
$ws.Range("K135").Value = 4257.75006
$ws.Range("L135").Value = 9330.3252
$ws.Range("M135").Value = -1722.75006
$ws.Range("N135").Value = -14400.3252

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2076.5
$ws.Range("I132").Value = 2025.5
$ws.Range("J132").Value = 2153
$ws.Range("K132").Value = 6076.5
$ws.Range("L132").Value = 6459
$ws.Range("M132").Value = -3546.5
$ws.Range("N132").Value = -11519

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3011.07
$ws.Range("I132").Value = 3151.796
$ws.Range("J132").Value = 2149.125
$ws.Range("K132").Value = 9455.387999999999
$ws.Range("L132").Value = 6447.375
$ws.Range("M132").Value = -6925.387999999999
$ws.Range("N132").Value = -11507.375

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3602.1897
$ws.Range("I132").Value = 4221.6978
$ws.Range("J132").Value = 1826.2667
$ws.Range("K132").Value = 12665.0934
$ws.Range("L132").Value = 5478.800099999999
$ws.Range("M132").Value = -10135.0934
$ws.Range("N132").Value = -10538.8001
